# Applies the cryptos list price/volume update described in the commit
# "Updated cryptos list on Fri Jun 16 20:10:21 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Row, $Col, $Text) {
    $cell = $ws.Cells.Item($Row, $Col)
    # Leading apostrophe forces Excel to treat numeric-looking strings
    # (e.g. "26.347.39", "1.000") as literal text instead of numbers/dates.
    $cell.Value = "'" + $Text
    # Reset the style so the quote-prefix flag added above does not linger
    # as a visible formatting change on the cell.
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "26.347.39"
Set-TextValue 2 5 "  +3.58%  "
Set-TextValue 3 4 "1.718.66"
Set-TextValue 3 5 "  +3.14%  "
Set-TextValue 4 4 "0.9995"
Set-TextValue 4 5 "  +0.11%  "
Set-TextValue 5 4 "239.53"
Set-TextValue 5 5 "  +1.14%  "
Set-TextValue 6 4 "1.000"
Set-TextValue 6 5 "  +0.08%  "
Set-TextValue 7 4 "0.4716"
Set-TextValue 7 5 "  -1.76%  "
Set-TextValue 8 4 "0.2634"
Set-TextValue 8 5 "  +0.69%  "
Set-TextValue 9 4 "0.06227"
Set-TextValue 9 5 "  +1.22%  "
Set-TextValue 10 4 "1.710.96"
Set-TextValue 10 5 "  +2.63%  "
Set-TextValue 11 4 "0.07082"
Set-TextValue 11 5 "  -0.02%  "
Set-TextValue 12 4 "15.30"
Set-TextValue 12 5 "  +3.31%  "
Set-TextValue 13 4 "0.5933"
Set-TextValue 13 5 "  -0.10%  "
Set-TextValue 14 4 "4.416"
Set-TextValue 14 5 "  +0.52%  "
Set-TextValue 15 4 "76.42"
Set-TextValue 15 5 "  +2.59%  "
Set-TextValue 16 4 "1.000"
Set-TextValue 16 5 "  +0.04%  "
Set-TextValue 17 4 "1.000"
Set-TextValue 17 5 "  +0.13%  "
Set-TextValue 18 4 "26.357.56"
Set-TextValue 18 5 "  +3.64%  "
Set-TextValue 19 4 "0.000006816"
Set-TextValue 19 5 "  +0.97%  "
Set-TextValue 20 4 "11.59"
Set-TextValue 20 5 "  +1.42%  "
Set-TextValue 21 4 "1.932.42"
Set-TextValue 21 5 "  +3.18%  "
Set-TextValue 22 4 "4.553"
Set-TextValue 22 5 "  +2.54%  "
Set-TextValue 23 4 "8.805"
Set-TextValue 23 5 "  +1.65%  "
Set-TextValue 24 4 "5.354"
Set-TextValue 24 5 "  +0.18%  "
Set-TextValue 25 4 "135.27"
Set-TextValue 25 5 "  +1.31%  "
Set-TextValue 26 5 "  +0.77%  "
Set-TextValue 27 5 "  +0.42%  "
Set-TextValue 28 4 "1.762"
Set-TextValue 28 5 "  +3.53%  "
Set-TextValue 29 4 "106.65"
Set-TextValue 29 5 "  +2.12%  "
Set-TextValue 30 4 "4.052"
Set-TextValue 30 5 "  +1.45%  "
Set-TextValue 31 4 "3.696"
Set-TextValue 31 5 "  +2.31%  "
Set-TextValue 32 4 "0.07725"
Set-TextValue 32 5 "  +0.96%  "
Set-TextValue 33 4 "0.04446"
Set-TextValue 33 5 "  +1.23%  "
Set-TextValue 34 4 "2.612"
Set-TextValue 34 5 "  +0.44%  "
Set-TextValue 35 4 "0.6218"
Set-TextValue 35 5 "  +2.03%  "
Set-TextValue 36 4 "0.9744"
Set-TextValue 36 5 "  +3.05%  "
Set-TextValue 37 4 "0.9320"
Set-TextValue 37 5 "  +8.61%  "
Set-TextValue 38 4 "115.23"
Set-TextValue 38 5 "  +16.90%  "
Set-TextValue 39 4 "2.411"
Set-TextValue 39 5 "  -8.22%  "
Set-TextValue 40 2 "RenderToken"
Set-TextValue 40 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 40 4 "1.910"
Set-TextValue 40 5 "  +4.50%  "
Set-TextValue 41 2 "PaxDollar"
Set-TextValue 41 3 "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue 41 4 "1.000"
Set-TextValue 41 5 "  +0.01%  "
Set-TextValue 42 4 "0.01473"
Set-TextValue 42 5 "  -1.97%  "
Set-TextValue 43 4 "5.280"
Set-TextValue 43 5 "  +12.68%  "
Set-TextValue 44 4 "0.3817"
Set-TextValue 44 5 "  +1.02%  "
Set-TextValue 45 5 "  +2.75%  "
Set-TextValue 46 5 "  +0.55%  "
Set-TextValue 47 4 "0.05290"
Set-TextValue 48 5 "  +3.24%  "
Set-TextValue 49 4 "7.698"
Set-TextValue 49 5 "  +4.71%  "
Set-TextValue 50 2 "Decentraland"
Set-TextValue 50 3 "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue 50 4 "0.3389"
Set-TextValue 50 5 "  +1.17%  "
Set-TextValue 51 2 "NEARProtocol"
Set-TextValue 51 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 51 4 "1.221"
Set-TextValue 51 5 "  +1.04%  "
